$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; existing rows 31-45 shift down to 32-46.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly price observation.
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value = "Bíobío"
$ws.Cells.Item(31, 4).Value = 44875
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100107
$ws.Cells.Item(31, 8).Value = "Otros"
$ws.Cells.Item(31, 9).Value = 100107002
$ws.Cells.Item(31, 10).Value = "Chirimoya"
$ws.Cells.Item(31, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 130
$ws.Cells.Item(31, 14).Value = 22000
$ws.Cells.Item(31, 15).Value = 23000
$ws.Cells.Item(31, 16).Value = 22538
$ws.Cells.Item(31, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(31, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 19).Value = 2254
$ws.Cells.Item(31, 20).Value = 10
